# Sprint 0 TODO list: break down the main activities into sub-tasks,
# so the parent tasks' estimated time drops to 0 (already covered by the
# new child rows), update column width, and fix up view/selection state.

$wb = $excel.ActiveWorkbook

$wsTRD = $wb.Worksheets.Item("TRD")
$wsDiseno = $wb.Worksheets.Item("Diseño")
$wsRecursos = $wb.Worksheets.Item("Recursos")

# --- TRD sheet: widen the "Descripcion" column (C) ---
$wsTRD.Columns.Item(3).ColumnWidth = 30.333333333333332

# --- TRD sheet: rows 4-24, add the new "A" (Tarea) / "B" (Funcion) values
#     and zero out the Tiempo Estimado (D) for the tasks that were broken
#     down into sub-tasks. ---

# Row 4
$wsTRD.Cells.Item(4, 1).Value = 3
$wsTRD.Cells.Item(4, 2).Value = 2
$wsTRD.Cells.Item(4, 4).Value = 0

# Row 5
$wsTRD.Cells.Item(5, 1).Value = 4
$wsTRD.Cells.Item(5, 2).Value = 3
$wsTRD.Cells.Item(5, 4).Value = 0

# Row 6
$wsTRD.Cells.Item(6, 1).Value = 5
$wsTRD.Cells.Item(6, 2).Value = 4
$wsTRD.Cells.Item(6, 4).Value = 0

# Row 7
$wsTRD.Cells.Item(7, 1).Value = 6
$wsTRD.Cells.Item(7, 2).Value = 5
$wsTRD.Cells.Item(7, 4).Value = 0

# Row 8
$wsTRD.Cells.Item(8, 1).Value = 7
$wsTRD.Cells.Item(8, 2).Value = 6
$wsTRD.Cells.Item(8, 4).Value = 0

# Row 9
$wsTRD.Cells.Item(9, 1).Value = 8
$wsTRD.Cells.Item(9, 2).Value = 1
$wsTRD.Cells.Item(9, 4).Value = 0

# Row 10 (new A/B only, D stays 0.5)
$wsTRD.Cells.Item(10, 1).Value = 9
$wsTRD.Cells.Item(10, 2).Value = 1

# Row 11
$wsTRD.Cells.Item(11, 1).Value = 10
$wsTRD.Cells.Item(11, 2).Value = 1

# Row 12
$wsTRD.Cells.Item(12, 1).Value = 11
$wsTRD.Cells.Item(12, 2).Value = 2

# Row 13
$wsTRD.Cells.Item(13, 1).Value = 12
$wsTRD.Cells.Item(13, 2).Value = 2

# Row 14
$wsTRD.Cells.Item(14, 1).Value = 13
$wsTRD.Cells.Item(14, 2).Value = 2

# Row 15
$wsTRD.Cells.Item(15, 1).Value = 14
$wsTRD.Cells.Item(15, 2).Value = 3

# Row 16
$wsTRD.Cells.Item(16, 1).Value = 15
$wsTRD.Cells.Item(16, 2).Value = 3

# Row 17
$wsTRD.Cells.Item(17, 1).Value = 16
$wsTRD.Cells.Item(17, 2).Value = 4

# Row 18
$wsTRD.Cells.Item(18, 1).Value = 17
$wsTRD.Cells.Item(18, 2).Value = 4

# Row 19
$wsTRD.Cells.Item(19, 1).Value = 18
$wsTRD.Cells.Item(19, 2).Value = 6

# Row 20
$wsTRD.Cells.Item(20, 1).Value = 19
$wsTRD.Cells.Item(20, 2).Value = 6

# Row 21
$wsTRD.Cells.Item(21, 1).Value = 20
$wsTRD.Cells.Item(21, 2).Value = 5

# Row 22
$wsTRD.Cells.Item(22, 1).Value = 21
$wsTRD.Cells.Item(22, 2).Value = 5

# Row 23
$wsTRD.Cells.Item(23, 1).Value = 22
$wsTRD.Cells.Item(23, 2).Value = 5

# Row 24
$wsTRD.Cells.Item(24, 1).Value = 23
$wsTRD.Cells.Item(24, 2).Value = 0

# --- Diseño sheet: move the selection to F8 (without activating the tab) ---
$wsDiseno.Range("F8").Select()

# --- TRD sheet: reset the view - drop the old scroll position / selection
#     and make TRD the active tab (instead of Recursos). ---
$wsTRD.Range("A1").Select()
$wsTRD.Activate()
